$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet from "stocks" to "stock" ---
$ws.Name = "stock"

# --- Preserve the existing header/index cell format (bold, bordered,
#     centered) before the old layout is wiped, by copying it to a
#     scratch cell well outside the range we are about to touch.
$ws.Range("B3").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Clear the old layout (old headers were on row 3, data rows 4-8,
#     spanning columns B:F) ---
$ws.Range("A1:F8").Clear()

# --- New layout: headers now on row 1, data on rows 2-6, with a new
#     0-based index column in A. ---

# Header row
$ws.Range("B1").Value = "tickers"
$ws.Range("C1").Value = "eps"
$ws.Range("D1").Value = "revenue"
$ws.Range("E1").Value = "price"
$ws.Range("F1").Value = "people"

# Row 2 - GOOGL
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "GOOGL"
$ws.Range("C2").Value = 27.82
$ws.Range("D2").Value = 87
$ws.Range("E2").Value = 845
$ws.Range("F2").Value = "larry page"

# Row 3 - WMT
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "WMT"
$ws.Range("C3").Value = 4.61
$ws.Range("D3").Value = 484
$ws.Range("E3").Value = 65
$ws.Range("F3").Value = "Allif"

# Row 4 - MSFT
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "MSFT"
$ws.Range("C4").Value = -1
$ws.Range("D4").Value = 85
$ws.Range("E4").Value = 64
$ws.Range("F4").Value = "bill gates"

# Row 5 - RIL (no eps value here)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "RIL "
$ws.Range("D5").Value = 50
$ws.Range("E5").Value = 1023
$ws.Range("F5").Value = "mukesh ambani"

# Row 6 - TATA (no price value here)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "TATA"
$ws.Range("C6").Value = 5.6
$ws.Range("D6").Value = -1
$ws.Range("F6").Value = "ratan tata"

# --- Re-apply the preserved header-style formatting to the new header
#     row and to the new index column (both use the bold/bordered style
#     that row 3 previously had). ---
$ws.Range("Z1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("A2:A6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Remove the scratch cell used to stash the format
$ws.Range("Z1").Clear()
